$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18 - this shifts the existing rows 18-38 down
# to 19-39 (and extends the used range to row 39), matching the diff's
# <dimension ref="A1:R39"/> and the row-by-row downward shift of all the
# historical price records.
$ws.Rows(18).Insert()

# Populate the newly inserted row 18 with the new weekly record.
$ws.Cells.Item(18, 1).Value = 4
$ws.Cells.Item(18, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(18, 3).Value = "Los Lagos"
$ws.Cells.Item(18, 4).Value = 44638
$ws.Cells.Item(18, 5).Value = 10
$ws.Cells.Item(18, 6).Value = 100112030
$ws.Cells.Item(18, 7).Value = "Poroto granado"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 60
$ws.Cells.Item(18, 11).Value = 27000
$ws.Cells.Item(18, 12).Value = 27000
$ws.Cells.Item(18, 13).Value = 27000
$ws.Cells.Item(18, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(18, 15).Value = "Región Metropolitana"
$ws.Cells.Item(18, 16).Value = 1080
$ws.Cells.Item(18, 17).Value = 25
$ws.Cells.Item(18, 18).Value = "Hortaliza"
